$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (PriceD / Volume(1h)E columns refreshed by the scraper)
$updates = @{
    'D2' = '30.501.14'
    'E2' = '  -0.49%  '
    'D3' = '1.885.38'
    'E3' = '  +0.65%  '
    'E4' = '  +0.15%  '
    'D5' = '243.38'
    'E5' = '  -1.86%  '
    'E6' = '  +0.17%  '
    'D7' = '0.4685'
    'E7' = '  -1.05%  '
    'D8' = '0.2896'
    'E8' = '  -0.69%  '
    'D9' = '0.06488'
    'E9' = '  +0.12%  '
    'D10' = '22.24'
    'E10' = '  +0.76%  '
    'D11' = '0.07746'
    'D12' = '1.887.61'
    'E12' = '  +0.83%  '
    'D13' = '95.40'
    'E13' = '  -1.29%  '
    'D14' = '0.7283'
    'E14' = '  -1.37%  '
    'D15' = '5.178'
    'E15' = '  +0.43%  '
    'D16' = '281.88'
    'E16' = '  +3.20%  '
    'D17' = '30.497.86'
    'E17' = '  -0.48%  '
    'D18' = '13.01'
    'E18' = '  -2.37%  '
    'D19' = '1.001'
    'E19' = '  +0.09%  '
    'D20' = '0.000007480'
    'E20' = '  -0.58%  '
    'D21' = '2.128.69'
    'E21' = '  +0.66%  '
    'D22' = '1.002'
    'E22' = '  +0.33%  '
    'D23' = '5.261'
    'E23' = '  -0.02%  '
    'D24' = '6.255'
    'E24' = '  +1.26%  '
    'D25' = '163.40'
    'E25' = '  -0.24%  '
    'D26' = '9.088'
    'E26' = '  -1.52%  '
    'E27' = '  +0.80%  '
    'D28' = '1.893'
    'E28' = '  -1.09%  '
    'D29' = '1.333'
    'E29' = '  -1.63%  '
    'D30' = '0.09711'
    'E30' = '  -2.79%  '
    'D31' = '1.469'
    'E31' = '  -2.77%  '
    'D32' = '4.280'
    'E32' = '  -0.05%  '
    'D33' = '4.131'
    'E33' = '  +0.67%  '
    'D34' = '0.04858'
    'E34' = '  +1.40%  '
    'E35' = '  +0.30%  '
    'D36' = '0.6934'
    'E36' = '  -0.45%  '
    'E37' = '  +0.17%  '
    'D38' = '0.01891'
    'E38' = '  +2.21%  '
    'D39' = '2.832'
    'E39' = '  +2.76%  '
    'D40' = '75.72'
    'E40' = '  +3.42%  '
    'D41' = '6.194'
    'E41' = '  -0.43%  '
    'D42' = '2.005'
    'E42' = '  +1.75%  '
    'D43' = '0.4253'
    'E43' = '  +1.61%  '
    'D44' = '1.001'
    'E44' = '  +0.11%  '
    'D45' = '0.8243'
    'E45' = '  -1.04%  '
    'D46' = '101.35'
    'E46' = '  -0.44%  '
    'D47' = '9.518'
    'E47' = '  +2.21%  '
    'D48' = '6.968'
    'E48' = '  -0.16%  '
    'D49' = '35.18'
    'E49' = '  -0.78%  '
    'D50' = '916.78'
    'E50' = '  -0.11%  '
    'D51' = '0.05752'
    'E51' = '  +1.83%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text semantics so numeric-looking strings (e.g. "1.001") are not
    # auto-coerced into Excel numbers, then drop the temporary format so the
    # cell keeps its original (default) style, matching the source data.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
